# Applies the "24th April 1st update" change:
#   - Inserts two new date columns (21/04/2020 and 22/04/2020) into the
#     long-format COVID state table, right after column AQ (21/03/2020),
#     shifting all the later date columns (old AR..BB) two places to the
#     right (new AS..BD).
#   - Populates the two new columns with the day's per-state counts.
#   - Fixes a data point in the existing 20/04/2020 column (AP16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the two new columns. Inserting at AR first shifts old AR.. right by
# one; inserting again at (the now second) AT shifts everything from there
# right by one more, leaving two fresh blank columns at AR and AT.
$ws.Columns("AR").Insert()
$ws.Columns("AT").Insert()

# Header row: label the two new columns with their dates (stored as plain
# text, matching how every other date header in row 1 is stored).
$ws.Range("AR1").Value = "21/04/2020"
$ws.Range("AT1").Value = "22/04/2020"

# New data for column AR (21/04/2020), keyed by row number.
$col_AR = @{
    3  = 1
    4  = 35
    7  = 13
    8  = 1
    10 = 75
    12 = 239
    13 = 4
    15 = 12
    17 = 10
    18 = 19
    20 = 67
    21 = 552
    23 = 1
    25 = 5
    27 = 6
    28 = 159
    29 = 76
    30 = 56
    32 = 153
    34 = 53
}

# New data for column AT (22/04/2020), keyed by row number.
$col_AT = @{
    25 = 3
    34 = 31
}

foreach ($row in $col_AR.Keys) {
    $ws.Range("AR$row").Value = $col_AR[$row]
}

foreach ($row in $col_AT.Keys) {
    $ws.Range("AT$row").Value = $col_AT[$row]
}

# Standalone correction to the existing 20/04/2020 figure for row 16.
$ws.Range("AP16").Value = 5
